$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldPrefix = "n425ac992496e4f4c872a808ed79f4ef9b"
$newPrefix = "n1e446dc28bab431fa3520a0ab74536d6b"

for ($row = 2; $row -le 22; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = [string]$cell.Value2
    if ($current.StartsWith($oldPrefix)) {
        $suffix = $current.Substring($oldPrefix.Length)
        $cell.Value2 = $newPrefix + $suffix
    }
}
